# Workbook state: a single "data" sheet with PanelApp gene rows.
# This change (per commit message "Refined metadata to be additional tab")
# adds a second "metadata" sheet summarizing the query, and refreshes the
# "time_taken" timestamps on the "data" sheet (column F) to the time the
# metadata tab's values were captured.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Refresh the "time_taken" (column F) timestamps on the "data" sheet.
# ---------------------------------------------------------------------
$timeTaken = @(
    "2021-10-05 14:21:46.283577",
    "2021-10-05 14:21:46.283585",
    "2021-10-05 14:21:46.283588",
    "2021-10-05 14:21:46.283590",
    "2021-10-05 14:21:46.283593",
    "2021-10-05 14:21:46.283596",
    "2021-10-05 14:21:46.283599",
    "2021-10-05 14:21:46.283601",
    "2021-10-05 14:21:46.283604",
    "2021-10-05 14:21:46.283606",
    "2021-10-05 14:21:46.283609",
    "2021-10-05 14:21:46.283611",
    "2021-10-05 14:21:46.283614",
    "2021-10-05 14:21:46.283616",
    "2021-10-05 14:21:46.283619",
    "2021-10-05 14:21:46.283621",
    "2021-10-05 14:21:46.283624",
    "2021-10-05 14:21:46.283626",
    "2021-10-05 14:21:46.283629",
    "2021-10-05 14:21:46.283631",
    "2021-10-05 14:21:46.283634",
    "2021-10-05 14:21:46.283636",
    "2021-10-05 14:21:46.283639",
    "2021-10-05 14:21:46.283641",
    "2021-10-05 14:21:46.283644",
    "2021-10-05 14:21:46.283646",
    "2021-10-05 14:21:46.283649",
    "2021-10-05 14:21:46.283651",
    "2021-10-05 14:21:46.283654",
    "2021-10-05 14:21:46.283656",
    "2021-10-05 14:21:46.283659",
    "2021-10-05 14:21:46.283661",
    "2021-10-05 14:21:46.283664",
    "2021-10-05 14:21:46.283667",
    "2021-10-05 14:21:46.283669",
    "2021-10-05 14:21:46.283672",
    "2021-10-05 14:21:46.283674",
    "2021-10-05 14:21:46.283676",
    "2021-10-05 14:21:46.283679",
    "2021-10-05 14:21:46.283681",
    "2021-10-05 14:21:46.283684",
    "2021-10-05 14:21:46.283687",
    "2021-10-05 14:21:46.283689",
    "2021-10-05 14:21:46.283692",
    "2021-10-05 14:21:46.283694",
    "2021-10-05 14:21:46.283696",
    "2021-10-05 14:21:46.283699",
    "2021-10-05 14:21:46.283701",
    "2021-10-05 14:21:46.283704",
    "2021-10-05 14:21:46.283706",
    "2021-10-05 14:21:46.283709",
    "2021-10-05 14:21:46.283711",
    "2021-10-05 14:21:46.283714",
    "2021-10-05 14:21:46.283717",
    "2021-10-05 14:21:46.283719",
    "2021-10-05 14:21:46.283721",
    "2021-10-05 14:21:46.283724",
    "2021-10-05 14:21:46.283726",
    "2021-10-05 14:21:46.283729",
    "2021-10-05 14:21:46.283731",
    "2021-10-05 14:21:46.283734",
    "2021-10-05 14:21:46.283736",
    "2021-10-05 14:21:46.283738",
    "2021-10-05 14:21:46.283741",
    "2021-10-05 14:21:46.283745"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 6).Value = $timeTaken[$i]
}

# ---------------------------------------------------------------------
# 2. Add a new "metadata" sheet after "data".
#    Duplicate "data" first so the new sheet inherits the same
#    sheet-level formatting (outline/page-setup/margins), then wipe its
#    contents and rebuild them from scratch.
# ---------------------------------------------------------------------
$ws1.Copy($null, $ws1) | Out-Null
$ws2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws2.Name = "metadata"
$ws2.Cells.Clear() | Out-Null

# Header row (bold / bordered style copied from the "data" sheet header).
$ws1.Range("B1").Copy() | Out-Null
$ws2.Range("B1:G1").PasteSpecial(-4122) | Out-Null

$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Row 2: the single metadata record, A2 styled like the "data" sheet's
# index column.
$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A2").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Value = 0

$ws2.Range("B2").Value = "Neurological ciliopathies"
$ws2.Range("C2").Value = 724

# "1.18" must stay literal text, not be coerced to the number 1.18.
# Compute it as a formula result in a scratch cell (T() forces a text
# result), then paste-special *values only* into D2 - this keeps the
# literal string without leaving D2's cell style altered.
$scratch = $ws2.Range("Z1")
$scratch.Formula = '=T("1.18")'
$scratch.Copy() | Out-Null
$ws2.Range("D2").PasteSpecial(-4163) | Out-Null
$scratch.Clear() | Out-Null

$ws2.Range("E2").Value = "2021-09-06T14:01:43.381315Z"
$ws2.Range("F2").Value = "2021-10-05 14:21:46.280163"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/724/?format=json"

# Leave the workbook's active sheet/selection as it was before ("data").
$ws1.Select() | Out-Null
$ws1.Range("A1").Select() | Out-Null
